# Updates cryptocurrency price ("Price", column D) and volume change
# ("Volume(1h)", column E) figures for the crypto list on the active
# worksheet, matching the latest scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.868.51"
$ws.Range("E2").Value = "  -4.17%  "
$ws.Range("D3").Value = "1.632.99"
$ws.Range("E3").Value = "  -6.47%  "
$ws.Range("D4").Value = "'0.9983"
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'233.11"
$ws.Range("E5").Value = "  -6.28%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.4715"
$ws.Range("E7").Value = "  -6.66%  "
$ws.Range("D8").Value = "'0.2553"
$ws.Range("E8").Value = "  -7.21%  "
$ws.Range("D9").Value = "'0.06067"
$ws.Range("E9").Value = "  -1.97%  "
$ws.Range("D10").Value = "'0.06958"
$ws.Range("E10").Value = "  -4.29%  "
$ws.Range("D11").Value = "1.641.99"
$ws.Range("E11").Value = "  -5.94%  "
$ws.Range("D12").Value = "'14.52"
$ws.Range("E12").Value = "  -4.43%  "
$ws.Range("D13").Value = "'0.5959"
$ws.Range("E13").Value = "  -9.04%  "
$ws.Range("D14").Value = "'4.305"
$ws.Range("E14").Value = "  -7.56%  "
$ws.Range("D15").Value = "'73.16"
$ws.Range("E15").Value = "  -5.92%  "
$ws.Range("D16").Value = "'0.9999"
$ws.Range("E16").Value = "  -0.01%  "
$ws.Range("D17").Value = "'0.9990"
$ws.Range("D18").Value = "24.874.73"
$ws.Range("E18").Value = "  -4.21%  "
$ws.Range("D19").Value = "'0.000006540"
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").Value = "'11.10"
$ws.Range("E20").Value = "  -6.45%  "
$ws.Range("D21").Value = "1.848.83"
$ws.Range("E21").Value = "  -6.08%  "
$ws.Range("D22").Value = "'4.321"
$ws.Range("E22").Value = "  -3.01%  "
$ws.Range("D23").Value = "'8.527"
$ws.Range("E23").Value = "  -2.31%  "
$ws.Range("E24").Value = "  -3.61%  "
$ws.Range("D25").Value = "'132.93"
$ws.Range("E25").Value = "  -2.77%  "
$ws.Range("E26").Value = "  -3.40%  "
$ws.Range("D27").Value = "'1.375"
$ws.Range("E27").Value = "  -8.63%  "
$ws.Range("D28").Value = "'103.32"
$ws.Range("E28").Value = "  -2.24%  "
$ws.Range("D29").Value = "'1.619"
$ws.Range("E29").Value = "  -9.24%  "
$ws.Range("D30").Value = "'3.805"
$ws.Range("E30").Value = "  -1.65%  "
$ws.Range("D31").Value = "'0.07695"
$ws.Range("E31").Value = "  -6.27%  "
$ws.Range("D32").Value = "'3.515"
$ws.Range("E32").Value = "  -3.77%  "
$ws.Range("D33").Value = "'0.9991"
$ws.Range("D34").Value = "'0.04275"
$ws.Range("E34").Value = "  -8.72%  "
$ws.Range("D35").Value = "'2.580"
$ws.Range("E35").Value = "  -2.83%  "
$ws.Range("D36").Value = "'0.9174"
$ws.Range("E36").Value = "  -7.97%  "
$ws.Range("D37").Value = "'0.5757"
$ws.Range("E37").Value = "  -7.04%  "
$ws.Range("D38").Value = "'2.531"
$ws.Range("E38").Value = "  -8.05%  "
$ws.Range("D39").Value = "'0.01533"
$ws.Range("E39").Value = "  -4.95%  "
$ws.Range("D40").Value = "'0.9984"
$ws.Range("E40").Value = "  -0.13%  "
$ws.Range("D41").Value = "'0.8043"
$ws.Range("E41").Value = "  +5.78%  "
$ws.Range("D42").Value = "'97.26"
$ws.Range("E42").Value = "  -3.47%  "
$ws.Range("D43").Value = "'1.759"
$ws.Range("E43").Value = "  -8.84%  "
$ws.Range("D44").Value = "'0.3670"
$ws.Range("E44").Value = "  -6.95%  "
$ws.Range("D45").Value = "'4.697"
$ws.Range("E45").Value = "  -6.31%  "
$ws.Range("D46").Value = "'0.05198"
$ws.Range("E46").Value = "  -1.60%  "
$ws.Range("D47").Value = "'0.1085"
$ws.Range("E47").Value = "  -5.69%  "
$ws.Range("D48").Value = "'5.983"
$ws.Range("E48").Value = "  -5.61%  "
$ws.Range("D49").Value = "'29.26"
$ws.Range("E49").Value = "  -4.69%  "
$ws.Range("D50").Value = "'0.9990"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "'0.9953"
$ws.Range("E51").Value = "  -0.45%  "
